$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "28.193.32"
$ws.Range("E2").Value = "  +2.24%  "
Set-TextValue "D3" "1.587.51"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +1.43%  "
Set-TextValue "D5" "213.14"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  +1.46%  "
Set-TextValue "D8" "23.92"
$ws.Range("E8").Value = "  +6.16%  "
$ws.Range("E9").Value = "  -0.13%  "
Set-TextValue "D10" "0.0598"
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("E11").Value = "  +2.58%  "
Set-TextValue "D12" "1.814.84"
$ws.Range("E12").Value = "  +1.28%  "
Set-TextValue "D13" "1.590.96"
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("E15").Value = "  -0.58%  "
Set-TextValue "D16" "28.230.87"
$ws.Range("E16").Value = "  +2.54%  "
Set-TextValue "D17" "63.10"
$ws.Range("E17").Value = "  +1.67%  "
Set-TextValue "D18" "227.13"
$ws.Range("E18").Value = "  +0.59%  "
Set-TextValue "D19" "0.0₃0707"
$ws.Range("E19").Value = "  +0.21%  "
Set-TextValue "D20" "7.46"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("E21").Value = "  +1.34%  "
Set-TextValue "D22" "4.09"
$ws.Range("E22").Value = "  -1.55%  "
Set-TextValue "D23" "9.31"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("E24").Value = "  +0.54%  "
Set-TextValue "D25" "151.79"
$ws.Range("E25").Value = "  +1.16%  "
Set-TextValue "D26" "15.16"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E27").Value = "  -0.60%  "
Set-TextValue "D28" "6.57"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("E29").Value = "  +1.39%  "
Set-TextValue "D30" "1.13"
$ws.Range("E30").Value = "  -0.27%  "
Set-TextValue "D31" "0.0471"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("E33").Value = "  -0.93%  "
Set-TextValue "D34" "1.395.95"
$ws.Range("E34").Value = "  -4.75%  "
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("E36").Value = "  -7.83%  "
$ws.Range("E37").Value = "  +1.65%  "
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").Value = "  +5.80%  "
Set-TextValue "D40" "0.541"
$ws.Range("E40").Value = "  -0.54%  "
Set-TextValue "D41" "0.811"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("E43").Value = "  +0.53%  "
Set-TextValue "D44" "5.60"
$ws.Range("E44").Value = "  -2.10%  "
Set-TextValue "D45" "0.979"
$ws.Range("E45").Value = "  +0.62%  "
Set-TextValue "D46" "64.21"
$ws.Range("E46").Value = "  -1.76%  "
Set-TextValue "D47" "1.724.21"
$ws.Range("E47").Value = "  +1.06%  "
$ws.Range("B48").Value = "mCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
Set-TextValue "D48" "2.13"
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D49" "87.13"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("E50").Value = "  +5.80%  "
$ws.Range("E51").Value = "  -0.39%  "
